# Austrian industrial statistics 1854 — apply the commit's edits:
#  1. Re-label the "l. Baumwoll - Spinnereien im Jahre 1854." block that
#     actually belongs to a Papier (paper) section (rows 58-120, col A) as
#     "m. Papier - Erzeugung im Jahre 1854." — this mints a new shared string.
#  2. Move the active tab from "note" back to "data", scroll the frozen
#     "data" view down toward the bottom of the sheet, and leave the last
#     used cell (B153) selected.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- 1. Fix the mislabeled rows: 58 through 120 in column A -----------------
$newLabel = "m. Papier - Erzeugung im Jahre 1854."
for ($r = 58; $r -le 120; $r++) {
    $data.Range("A" + $r).Value = $newLabel
}

# --- 2. Switch the active sheet / view state --------------------------------
# Select the top-left cell first (matches the pre-edit "B1" top-pane
# selection) before (re)applying the freeze so the frozen split stays at row 1.
$data.Activate()
[void]$data.Range("B1").Select()
$excel.ActiveWindow.FreezePanes = $true

# Scroll the bottom pane down near the end of the data and land the
# selection on B153, the last populated cell.
$excel.ActiveWindow.ScrollRow = 119
[void]$data.Range("A119:B153").Select()
[void]$data.Range("B153").Activate()

# "note" was the active tab before the edit (tabSelected=1); activating
# "data" above already flips tabSelected onto "data" and off of "note".
